$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.229.39"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "3.532.37"
$ws.Range("E3").Value = "  +2.77%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'597.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").Value = "'138.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "

$ws.Range("D7").Value = "3.530.72"
$ws.Range("E7").Value = "  +2.80%  "

$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  +2.59%  "

$ws.Range("D11").Value = "'6.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.60%  "

$ws.Range("D12").Value = "'0.388"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.69%  "

$ws.Range("D13").Value = "4.131.94"
$ws.Range("E13").Value = "  +2.57%  "

$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("D15").Value = "'27.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.73%  "

$ws.Range("D16").Value = "3.528.01"
$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("D18").Value = "65.275.08"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").Value = "'10.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.14%  "

$ws.Range("D20").Value = "'5.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "

$ws.Range("D21").Value = "'14.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.38%  "

$ws.Range("D22").Value = "'393.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'0.572"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.90%  "

$ws.Range("D24").Value = "3.671.87"
$ws.Range("E24").Value = "  +2.64%  "

$ws.Range("D25").Value = "'73.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  +7.17%  "

$ws.Range("D28").Value = "'7.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.90%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("E30").Value = "  +2.68%  "

$ws.Range("D31").Value = "'8.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("D32").Value = "3.545.03"
$ws.Range("E32").Value = "  +2.87%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "'23.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.44%  "

$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").Value = "'1.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.34%  "

$ws.Range("D37").Value = "'6.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.79%  "

$ws.Range("D38").Value = "'168.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.35%  "

$ws.Range("E39").Value = "  +4.50%  "

$ws.Range("D40").Value = "'4.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.02%  "

$ws.Range("D41").Value = "'0.0799"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.58%  "

$ws.Range("D42").Value = "'0.823"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").Value = "'25.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.13%  "

$ws.Range("D44").Value = "'42.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.17%  "

$ws.Range("D45").Value = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").Value = "'4.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("E47").Value = "  +3.40%  "

$ws.Range("D48").Value = "'1.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.84%  "

$ws.Range("D49").Value = "'6.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.88%  "

$ws.Range("D50").Value = "2.381.33"
$ws.Range("E50").Value = "  +7.68%  "

$ws.Range("D51").Value = "'302.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.50%  "

